$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 112178515
$ws.Range("B18").Value = 90480
$ws.Range("E18").Value = 4769
$ws.Range("F18").Value = "Svavelriska"
$ws.Range("G18").Value = "Lactarius scrobiculatus"
$ws.Range("H18").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q18").Value = 760089
$ws.Range("R18").Value = 7210467

$ws.Range("A19").Value = 112178516
$ws.Range("B19").Value = 89331
$ws.Range("E19").Value = 3215
$ws.Range("F19").Value = "Rödgul trumpetsvamp"
$ws.Range("G19").Value = "Craterellus lutescens"
$ws.Range("H19").Value = "(Fr.) Fr."
$ws.Range("Q19").Value = 760126
$ws.Range("R19").Value = 7210471

$ws.Range("A20").Value = 112178517
$ws.Range("B20").Value = 90480
$ws.Range("E20").Value = 4769
$ws.Range("F20").Value = "Svavelriska"
$ws.Range("G20").Value = "Lactarius scrobiculatus"
$ws.Range("H20").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q20").Value = 760128
$ws.Range("R20").Value = 7210459

$ws.Range("A21").Value = 112178538
$ws.Range("B21").Value = 98891
$ws.Range("E21").Value = 222771
$ws.Range("F21").Value = "Svart trolldruva"
$ws.Range("G21").Value = "Actaea spicata"
$ws.Range("H21").Value = "L."
$ws.Range("Q21").Value = 760363
$ws.Range("R21").Value = 7210127
